$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 592921.75
$ws.Range("J17").Value = 592921.75
$ws.Range("L17").Value = 1778765.25
$ws.Range("N17").Value = -1779101.25
$ws.Range("H93").Value = 11272.728
$ws.Range("J93").Value = 11272.728
$ws.Range("L93").Value = 11272.728
$ws.Range("N93").Value = -16264.728
$ws.Range("H116").Value = 4244.25
$ws.Range("I116").Value = 4816.364
$ws.Range("J116").Value = 3760.1538
$ws.Range("K116").Value = 4816.364
$ws.Range("L116").Value = 3760.1538
$ws.Range("M116").Value = -1374.364
$ws.Range("N116").Value = -10644.1538
$ws.Range("H132").Value = 2122.7144
$ws.Range("I132").Value = 2572.5
$ws.Range("K132").Value = 7717.5
$ws.Range("M132").Value = -5187.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4152.68
$ws.Range("I32").Value = 3375.0557
$ws.Range("J32").Value = 11151.3
$ws.Range("K32").Value = 3375.0557
$ws.Range("L32").Value = 11151.3
$ws.Range("M32").Value = -3088.0557
$ws.Range("N32").Value = -11725.3
$ws.Range("H74").Value = 992.32434
$ws.Range("I74").Value = 916.64703
$ws.Range("J74").Value = 1850
$ws.Range("K74").Value = 916.64703
$ws.Range("L74").Value = 1850
$ws.Range("M74").Value = -42.64702999999997
$ws.Range("N74").Value = -3598
$ws.Range("H77").Value = 992.32434
$ws.Range("I77").Value = 916.64703
$ws.Range("J77").Value = 1850
$ws.Range("K77").Value = 4583.23515
$ws.Range("L77").Value = 9250
$ws.Range("M77").Value = -215.2351499999995
$ws.Range("N77").Value = -17986
$ws.Range("H102").Value = 2112.577
$ws.Range("I102").Value = 2166.913
$ws.Range("J102").Value = 1696
$ws.Range("K102").Value = 2166.913
$ws.Range("L102").Value = 1696
$ws.Range("M102").Value = -544.913
$ws.Range("N102").Value = -4940
$ws.Range("H132").Value = 2022.5588
$ws.Range("I132").Value = 1422.1765
$ws.Range("J132").Value = 2622.9412
$ws.Range("K132").Value = 4266.529500000001
$ws.Range("L132").Value = 7868.823600000001
$ws.Range("M132").Value = -1736.529500000001
$ws.Range("N132").Value = -12928.8236

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 7322.3335
$ws.Range("I36").Value = 2433.3333
$ws.Range("J36").Value = 8952
$ws.Range("K36").Value = 2433.3333
$ws.Range("L36").Value = 8952
$ws.Range("M36").Value = -1899.3333
$ws.Range("N36").Value = -10020
$ws.Range("H37").Value = 10028.637
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 10028.637
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 10028.637
$ws.Range("M37").Value = ""
$ws.Range("N37").Value = -10302.637
$ws.Range("H39").Value = 9289.75
$ws.Range("J39").Value = 9289.75
$ws.Range("L39").Value = 9289.75
$ws.Range("N39").Value = -10067.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 1750
$ws.Range("J15").Value = 1750
$ws.Range("L15").Value = 1750
$ws.Range("N15").Value = -2090

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 1603.8572
$ws.Range("J64").Value = 1712.8334
$ws.Range("L64").Value = 5138.5002
$ws.Range("N64").Value = -5678.5002
$ws.Range("H67").Value = 1603.8572
$ws.Range("J67").Value = 1712.8334
$ws.Range("L67").Value = 5138.5002
$ws.Range("N67").Value = -7010.5002
$ws.Range("H96").Value = 4330.8945
$ws.Range("J96").Value = 4330.8945
$ws.Range("L96").Value = 12992.6835
$ws.Range("N96").Value = -17110.6835
$ws.Range("H129").Value = 35189.8
$ws.Range("I129").Value = 1062.4166
$ws.Range("J129").Value = 57941.39
$ws.Range("K129").Value = 3187.2498
$ws.Range("L129").Value = 173824.17
$ws.Range("M129").Value = 1812.7502
$ws.Range("N129").Value = -183824.17
$ws.Range("H131").Value = 7172291.5
$ws.Range("I131").Value = 83500440
$ws.Range("J131").Value = 16527.594
$ws.Range("K131").Value = 250501320
$ws.Range("L131").Value = 49582.78200000001
$ws.Range("M131").Value = -250496280
$ws.Range("N131").Value = -59662.78200000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 44831.332
$ws.Range("J123").Value = 44831.332
$ws.Range("L123").Value = 44831.332
$ws.Range("N123").Value = -49731.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 10000000
$ws.Range("I3").Value = 10000000
$ws.Range("K3").Value = 10000000
$ws.Range("M3").Value = -9999888
$ws.Range("H15").Value = 10000000
$ws.Range("I15").Value = 10000000
$ws.Range("K15").Value = 10000000
$ws.Range("M15").Value = -9999830
$ws.Range("H22").Value = 606.75
$ws.Range("I22").Value = 399.66666
$ws.Range("J22").Value = 731
$ws.Range("K22").Value = 399.66666
$ws.Range("L22").Value = 731
$ws.Range("M22").Value = -104.66666
$ws.Range("N22").Value = -1321
$ws.Range("H27").Value = 606.75
$ws.Range("I27").Value = 399.66666
$ws.Range("J27").Value = 731
$ws.Range("K27").Value = 399.66666
$ws.Range("L27").Value = 731
$ws.Range("M27").Value = -292.66666
$ws.Range("N27").Value = -945
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").Value = ""
$ws.Range("H122").Value = 2620.8
$ws.Range("I122").Value = 2026
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 6078
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -3628
$ws.Range("N122").Value = -19900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 20933.334
$ws.Range("J32").Value = 20933.334
$ws.Range("L32").Value = 20933.334
$ws.Range("N32").Value = -21567.334
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").Value = ""
$ws.Range("H107").Value = 1535.4736
$ws.Range("I107").Value = 1802.1818
$ws.Range("J107").Value = 1168.75
$ws.Range("K107").Value = 5406.5454
$ws.Range("L107").Value = 3506.25
$ws.Range("M107").Value = -3486.5454
$ws.Range("N107").Value = -7346.25
$ws.Range("H122").Value = 1399
$ws.Range("I122").Value = 1298.8
$ws.Range("J122").Value = 1900
$ws.Range("K122").Value = 3896.4
$ws.Range("L122").Value = 5700
$ws.Range("M122").Value = -1446.4
$ws.Range("N122").Value = -10600
